$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.470.95"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "3.665.00"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.73"
$ws.Range("E5").Value = "  -7.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.57"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.20"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  -2.87%  "
$ws.Range("D13").Value = "4.288.43"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.23"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "3.667.76"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "69.480.31"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.117"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.50"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.87"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.28"
$ws.Range("E20").Value = "  +4.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "468.08"
$ws.Range("E21").Value = "  -1.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.648"
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "79.63"
$ws.Range("E23").Value = "  -1.00%  "
$ws.Range("D24").Value = "3.815.07"
$ws.Range("E24").Value = "  -1.05%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.16"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  -5.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.61"
$ws.Range("E28").Value = "  -6.22%  "
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.57"
$ws.Range("E33").Value = "  -1.45%  "
$ws.Range("D34").Value = "3.673.45"
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "178.24"
$ws.Range("E39").Value = "  +2.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("E42").Value = "  -5.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0891"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.922"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.87"
$ws.Range("E45").Value = "  -0.51%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.19"
$ws.Range("E46").Value = "  +5.19%  "
$ws.Range("E47").Value = "  -2.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.85"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000264"
$ws.Range("E49").Value = "  -5.95%  "
$ws.Range("E50").Value = "  -5.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.259"
$ws.Range("E51").Value = "  -3.26%  "
